# Rename the "wt" and "dcin5" sheets to include the "_log2_expression" suffix,
# matching the commit message:
#   "Updated the input files in sixteen_tests to have strain_log2_expression
#    instead of just strain"
$wb = $excel.ActiveWorkbook

$wsWt = $wb.Worksheets.Item("wt")
$wsWt.Name = "wt_log2_expression"

$wsDcin5 = $wb.Worksheets.Item("dcin5")
$wsDcin5.Name = "dcin5_log2_expression"

# The diff also shows the active/selected sheet moving from
# "optimization_parameters" (old activeTab=6 / tabSelected on that sheet)
# to "dcin5" (new activeTab=3 / tabSelected="1" on the dcin5 sheet).
# Selecting the renamed dcin5 sheet reproduces that: it becomes the active
# sheet (updates workbookView's activeTab) and moves tabSelected="1" onto
# its sheetView while removing it from whichever sheet had it before.
$wsDcin5.Select()
